# Weekly refresh of the "Fruta, Vega Central Mapocho de Santiago - Coco" sheet.
# Each data row (2-41) gets its Fecha/Calidad/Volumen/Precio.../Origen/Precio-$-Kg
# block replaced by the corresponding block from another row of the same
# (pre-edit) sheet, per the mapping below (row -> source row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("D", "L", "M", "N", "O", "P", "R", "S")

# Snapshot the current (pre-edit) values for every touched column/row BEFORE
# writing anything, since several rows source their new values from rows
# that are themselves being overwritten.
$snapshot = @{}
for ($r = 2; $r -le 41; $r++) {
    $rowData = @{}
    foreach ($c in $cols) {
        $rowData[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $rowData
}

# Destination row -> source row (values copied from source's pre-edit state).
$rowMap = @{
    2 = 20
    3 = 2
    4 = 4
    5 = 24
    6 = 17
    7 = 30
    8 = 25
    9 = 28
    10 = 36
    11 = 12
    12 = 10
    13 = 38
    14 = 32
    15 = 22
    16 = 23
    17 = 31
    18 = 16
    19 = 13
    20 = 9
    21 = 18
    22 = 41
    23 = 6
    24 = 29
    25 = 27
    26 = 34
    27 = 26
    28 = 5
    29 = 39
    30 = 35
    31 = 14
    32 = 37
    33 = 21
    34 = 7
    35 = 8
    36 = 33
    37 = 11
    38 = 15
    39 = 3
    40 = 40
    41 = 19
}

foreach ($destRow in $rowMap.Keys) {
    $srcRow = $rowMap[$destRow]
    $srcData = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $ws.Range("$c$destRow").Value = $srcData[$c]
    }
}
